# Sprint 4 Backlog - Burndown: add planned meals time entry (row 6) that
# was missing since it wasn't brought in by the merge.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 corresponds to the "Add planned meal ingredients to shopping list"
# task which was missing its actual-time tracking values.
$ws.Range("E6").Value = 1.5
$ws.Range("F6").Value = $ws.Range("D6").Value2
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# Reflect the new selection state left after the edit.
$ws.Range("G6").Select() | Out-Null
